# Update main GSC export data: the oldest day (2025-11-13) has dropped out of
# the rolling export window, so the whole "Chart" sheet's data table shifts
# up by one row (each remaining day's row moves from row N to row N-1), and
# the newest day's row (old row 88 / 2026-02-07) now lands on row 87.
#
# Deleting the entire sheet row (rather than rewriting cell-by-cell) lets
# Excel itself shift every subsequent row's cells + styles up by one and
# keeps xl/sharedStrings.xml consistent (including dropping the
# now-duplicate "2025-11-14" shared string).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
